# corrigindo o saldo dos produtos
# Adds two new "ENTRADA" (stock-in) movement rows to the "movimentos" sheet,
# mirroring the formatting of the existing last row (row 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("movimentos")

function Set-TextCell($rng, [string]$text) {
    # Force the cell to hold a literal text value (not a number/date), while
    # keeping whatever cell style/format is already applied to the cell.
    $escaped = $text -replace '"', '""'
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)  # xlPasteValues
}

# Insert two new rows at the bottom (after row 19), copying the formatting
# of row 19 so the new rows keep the same cell style (s="8") across every
# column, including A and C which have no column-level default style.
$ws.Rows(19).Copy()
$ws.Rows(20).Insert(-4121)  # xlShiftDown, carries the copied formatting in
$ws.Rows(19).Copy()
$ws.Rows(21).Insert(-4121)
$excel.CutCopyMode = 0

# Row 20
$ws.Range("A20").Value = 19
Set-TextCell $ws.Range("B20") "333"
Set-TextCell $ws.Range("C20") "ENTRADA"
$ws.Range("D20").Value = 33
Set-TextCell $ws.Range("E20") "2026-01-20 11:21:20"

# Row 21
$ws.Range("A21").Value = 20
Set-TextCell $ws.Range("B21") "333"
Set-TextCell $ws.Range("C21") "ENTRADA"
$ws.Range("D21").Value = 3
Set-TextCell $ws.Range("E21") "2026-01-20 11:30:53"

$excel.CutCopyMode = 0
